$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New localization rows to append (Name, zh-CN, en-US)
$rows = @(
    @("feedback.send", "提交", "Submit"),
    @("feedback.hero", "我们会阅读每一条反馈", "We Read Every Feedback"),
    @("feedback.debug", "发送错误日志", "Send Error Log"),
    @("feedback.content.placeholder", "输入内容", "Write something you want to tell us"),
    @("feedback.contact.placeholder", "适合我们联系你的方式", "How can we contact you?")
)

$startRow = 28
$lastExistingRow = 27

# Row 32 is brand new - clone formatting (styles + row height) from the last
# existing data row (27) before writing values into it.
$newRowIndex = $startRow + $rows.Count - 1
$ws.Range("A$lastExistingRow`:E$lastExistingRow").Copy()
$ws.Range("A$newRowIndex`:E$newRowIndex").PasteSpecial(-4122)
$ws.Rows.Item($newRowIndex).RowHeight = 20.1

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $triple = $rows[$i]

    # Rows 28-31 already existed with the "empty data cell" style (s=4) on
    # A:C; switch them to the "filled data cell" style (s=3) used elsewhere
    # in the sheet by pasting the formatting from a populated row first.
    if ($r -ne $newRowIndex) {
        $ws.Range("A$lastExistingRow`:C$lastExistingRow").Copy()
        $ws.Range("A$r`:C$r").PasteSpecial(-4122)
    }

    $ws.Cells.Item($r, 1).Value = $triple[0]
    $ws.Cells.Item($r, 2).Value = $triple[1]
    $ws.Cells.Item($r, 3).Value = $triple[2]
}
